$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet
$ws.Name = "Лист"

# 2. Add the 2023 column (Q)
$ws.Range("Q4").Value = 2023
$ws.Range("Q4").Font.Bold = $true
$ws.Range("Q4").Font.Italic = $false
$ws.Range("Q4").Font.Name = "Times New Roman"
$ws.Range("Q4").Font.Size = 9
$ws.Range("Q4").HorizontalAlignment = -4152
$ws.Range("Q4").VerticalAlignment = -4108
$ws.Range("Q4").WrapText = $true
$ws.Range("Q4").Borders.Item(9).LineStyle = -4138
$ws.Range("Q4").Borders.Item(9).Weight = -4138
$ws.Range("Q4").Borders.Item(10).LineStyle = -4138
$ws.Range("Q4").Borders.Item(10).Weight = -4138

$ws.Range("Q5").Value = 3.2
$ws.Range("Q5").NumberFormat = "0.0"
$ws.Range("Q5").Font.Name = "Times New Roman"
$ws.Range("Q5").Font.Size = 9
$ws.Range("Q5").VerticalAlignment = -4108
$ws.Range("Q5").Borders.Item(9).LineStyle = -4138
$ws.Range("Q5").Borders.Item(9).Weight = -4138
$ws.Range("Q5").Borders.Item(10).LineStyle = -4138
$ws.Range("Q5").Borders.Item(10).Weight = -4138

# 3. Update row 4 headers (A4:C4) - new "Name of indicators" labels
$ws.Range("A4").Value = "Көрсөткүчтөрдүн аталышы"
$ws.Range("B4").Value = "Наименование показателей"
$ws.Range("C4").Value = "Items"

$ws.Range("A4:C4").Font.Bold = $true
$ws.Range("A4:C4").Font.Italic = $false
$ws.Range("A4:C4").Font.Name = "Times New Roman"
$ws.Range("A4:C4").Font.Size = 9
$ws.Range("A4:C4").HorizontalAlignment = -4108
$ws.Range("A4:C4").VerticalAlignment = -4108

# 4. Row 5: remove right alignment; add custom height
$ws.Rows("5").RowHeight = 40.5
$ws.Range("A5:Q5").VerticalAlignment = -4108

# 5. Row 6: footnotes with superscript "1" lead-in (rich text)
$ws.Rows("6").RowHeight = 37.5

$txtA = "1 2019-жылдан баштап маалыматтар, 2008 жылдагы Улуттук Эсептер Тутумунун эл аралык стандарттарына ылайык эсептелген "
$ws.Range("A6").Value = $txtA
$a6c1 = $ws.Range("A6").Characters(1, 1)
$a6c1.Font.Superscript = $true
$a6c1.Font.Italic = $true
$a6c1.Font.Bold = $false
$a6c1.Font.Size = 8
$a6c1.Font.Name = "Times New Roman"
$a6c2 = $ws.Range("A6").Characters(2, $txtA.Length - 1)
$a6c2.Font.Italic = $true
$a6c2.Font.Bold = $false
$a6c2.Font.Size = 8
$a6c2.Font.Name = "Times New Roman"
$ws.Range("A6").HorizontalAlignment = -4131
$ws.Range("A6").VerticalAlignment = -4108
$ws.Range("A6").WrapText = $true

$txtB = " Данные с 2019 года рассчитаны по международному стандарту Системы Национальных Счетов 2008 года"
$ws.Range("B6").Value = "1" + $txtB
$b6c1 = $ws.Range("B6").Characters(1, 1)
$b6c1.Font.Superscript = $true
$b6c1.Font.Italic = $true
$b6c1.Font.Bold = $false
$b6c1.Font.Size = 8
$b6c1.Font.Name = "Times New Roman"
$b6c2 = $ws.Range("B6").Characters(2, $txtB.Length)
$b6c2.Font.Italic = $true
$b6c2.Font.Bold = $false
$b6c2.Font.Size = 8
$b6c2.Font.Name = "Times New Roman"
$ws.Range("B6").HorizontalAlignment = -4131
$ws.Range("B6").WrapText = $true

$txtC = "Data from 2019 are calculated according to the international standard of the System of National Accounts 2008"
$ws.Range("C6").Value = "1 " + $txtC
$c6c1 = $ws.Range("C6").Characters(1, 2)
$c6c1.Font.Superscript = $true
$c6c1.Font.Italic = $true
$c6c1.Font.Bold = $false
$c6c1.Font.Size = 8
$c6c1.Font.Name = "Times New Roman"
$c6c2 = $ws.Range("C6").Characters(3, $txtC.Length)
$c6c2.Font.Italic = $true
$c6c2.Font.Bold = $false
$c6c2.Font.Size = 8
$c6c2.Font.Name = "Times New Roman"
$ws.Range("C6").HorizontalAlignment = -4131
$ws.Range("C6").VerticalAlignment = -4108
$ws.Range("C6").WrapText = $true

# 6. Update selection to match target
$ws.Range("Q9").Select()
